# Apply the "version final sin errores" edit:
#  - Metadata sheet: bump Version, Date, and Context values.
#  - Elements sheet: the root Extension row (row 1) now also carries the
#    ele-1/ext-1 invariant text (already present on the Extension.extension
#    row, AJ3) in its Invariants column (AJ1).

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "0.7.0"
$meta.Range("B8").Value = "2023-09-13T17:11:14-03:00"
$meta.Range("B20").Value = "element:Patient"

$elements = $wb.Worksheets.Item("Elements")
$invariantText = $elements.Range("AJ3").Value()
$elements.Range("AJ1").Value = $invariantText
